$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 over to I1:J1,
# then set the new header labels.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 / IF numeric columns for rows 2-14.
$values = @(
    @(5, 6),
    @(1, 4),
    @(2, 5),
    @(8, 8),
    @(3, 5),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(4, 6),
    @(5, 7)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
